# Update the US scaling mapping workbook:
#  - Rename several scaling_sector / ceds_sector identifiers to use
#    hyphens instead of underscores in certain tokens (industry_comb ->
#    industry-comb, etc.)
#  - Split out "off-highway" and "other-end-use-sectors" scaling
#    sectors from what used to be lumped under "FUEL COMB. OTHER" /
#    "OFF-HIGHWAY" / "industry_comb", adding new ceds_sector mapping
#    rows (43-49) on the "map" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# --- Top section renames / additions ---
$ws.Cells.Item(2, 2).Value2 = "industry-comb"
$ws.Cells.Item(3, 2).Value2 = "off-highway"
$ws.Cells.Item(5, 2).Value2 = "1A1a-Electricity-public"
$ws.Cells.Item(5, 3).Value2 = "1A1a_Electricity-public"
$ws.Cells.Item(6, 2).Value2 = "1A3b-Road"
$ws.Cells.Item(7, 2).Value2 = "other-end-use-sectors"
$ws.Cells.Item(9, 2).Value2 = "other-industrial-process"
$ws.Cells.Item(10, 2).Value2 = "1A1bc-Other-transformation"
$ws.Cells.Item(12, 2).Value2 = "2C-Metal-production"
$ws.Cells.Item(20, 2).Value2 = "2B-Chemical-industry"

# --- industry_comb -> industry-comb rows ---
$ws.Cells.Item(23, 2).Value2 = "industry-comb"
$ws.Cells.Item(24, 2).Value2 = "industry-comb"
$ws.Cells.Item(25, 2).Value2 = "industry-comb"
$ws.Cells.Item(26, 2).Value2 = "industry-comb"
$ws.Cells.Item(27, 2).Value2 = "industry-comb"
$ws.Cells.Item(28, 2).Value2 = "industry-comb"
$ws.Cells.Item(29, 2).Value2 = "industry-comb"

# Row 29 (was 1A2g_Ind-Comb-Construction) drops out of industry-comb and
# the remaining ceds_sector values shift up by one row; the construction
# entry reappears later under off-highway (row 43).
$ws.Cells.Item(29, 3).Value2 = "1A2g_Ind-Comb-transpequip"
$ws.Cells.Item(30, 2).Value2 = "industry-comb"
$ws.Cells.Item(30, 3).Value2 = "1A2g_Ind-Comb-machinery"
$ws.Cells.Item(31, 2).Value2 = "industry-comb"
$ws.Cells.Item(31, 3).Value2 = "1A2g_Ind-Comb-mining-quarying"
$ws.Cells.Item(32, 2).Value2 = "industry-comb"
$ws.Cells.Item(32, 3).Value2 = "1A2g_Ind-Comb-wood-products"
$ws.Cells.Item(33, 2).Value2 = "industry-comb"
$ws.Cells.Item(33, 3).Value2 = "1A2g_Ind-Comb-textile-leather"
$ws.Cells.Item(34, 2).Value2 = "industry-comb"
$ws.Cells.Item(34, 3).Value2 = "1A2g_Ind-Comb-other"

$ws.Cells.Item(35, 2).Value2 = "waste"
$ws.Cells.Item(35, 3).Value2 = "5A_Solid-waste-disposal"
$ws.Cells.Item(36, 3).Value2 = "5E_Other-waste-handling"
$ws.Cells.Item(37, 3).Value2 = "5C_Waste-incineration"

$ws.Cells.Item(38, 2).Value2 = "solvents"
$ws.Cells.Item(38, 3).Value2 = "2D3_Other-product-use"
$ws.Cells.Item(39, 3).Value2 = "2D_Paint-application"
$ws.Cells.Item(40, 3).Value2 = "2D_Degreasing-Cleaning"
$ws.Cells.Item(41, 3).Value2 = "2D3_Chemical-product-use"
$ws.Cells.Item(42, 3).Value2 = "2D3_Other-product-use"

# --- New rows 43-49: off-highway and other-end-use-sectors detail ---
$ws.Cells.Item(43, 2).Value2 = "off-highway"
$ws.Cells.Item(43, 3).Value2 = "1A2g_Ind-Comb-Construction"
$ws.Cells.Item(44, 2).Value2 = "off-highway"
$ws.Cells.Item(44, 3).Value2 = "1A3c_Rail"
$ws.Cells.Item(45, 2).Value2 = "off-highway"
$ws.Cells.Item(45, 3).Value2 = "1A3eii_Other-transp"
$ws.Cells.Item(46, 2).Value2 = "off-highway"
$ws.Cells.Item(46, 3).Value2 = "1A4c_Agriculture-forestry-fishing"

$ws.Cells.Item(47, 2).Value2 = "other-end-use-sectors"
$ws.Cells.Item(47, 3).Value2 = "1A4a_Commercial-institutional"
$ws.Cells.Item(48, 2).Value2 = "other-end-use-sectors"
$ws.Cells.Item(48, 3).Value2 = "1A4b_Residential"
$ws.Cells.Item(49, 2).Value2 = "other-end-use-sectors"
$ws.Cells.Item(49, 3).Value2 = "1A5_Other-unspecified"

# Match the author's final selection position on the map sheet.
$ws.Range("C49").Select()
